# incorporate SDRS policy into DB
# Applies the edits described by the commit "incorporate SDRS policy into DB":
#  - turns on (TRUE) the EEC_return / EEC_FR / EEC_sharedADC corrective-action
#    rows (B10:B12) on the paramlist sheet
#  - adds a new DC_hybrid policy row (row 16), cloned from the existing
#    cola_SDRS row (row 13) but switched to a constant COLA / DC-style setup
#  - adds a new DC_pure policy row (row 17) with just a name in column A
#  - extends the H and F column dropdown (data validation) ranges to include
#    the new row 16
#  - updates the active selection to reflect the new edit location

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("paramlist")

# --- Flip the three corrective-action switches from FALSE to TRUE ---
$ws.Range("B10").Value = $true
$ws.Range("B11").Value = $true
$ws.Range("B12").Value = $true

# --- New row 16: DC_hybrid (clone of row 13 "cola_SDRS" with tweaks) ---
$ws.Range("A16").Value = "DC_hybrid"
$ws.Range("B16").Value = $false
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = "constant"
$ws.Range("F16").Value = "fixed"
$ws.Range("G16").Value = 0.015
$ws.Range("H16").Value = "ALpct"
$ws.Range("I16").Value = 1
$ws.Range("J16").Value = 1
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.06
$ws.Range("M16").Value = 0.06
$ws.Range("R16").Value = 0.02
$ws.Range("S16").Value = 0
$ws.Range("T16").Value = 0.02
$ws.Range("U16").Value = 0
$ws.Range("V16").Value = 0.9
$ws.Range("W16").Value = 0.02
$ws.Range("X16").Value = 0
$ws.Range("Y16").Value = 0.9
$ws.Range("Z16").Value = 0.01
$ws.Range("AA16").Value = 0.001
$ws.Range("AB16").Value = 0.02
$ws.Range("AC16").Value = 0
$ws.Range("AD16").Value = 1
$ws.Range("AE16").Formula = "=0.15/20"
$ws.Range("AF16").Value = 0.001
$ws.Range("AG16").Value = 0.02
$ws.Range("AH16").Value = 0
$ws.Range("AI16").Value = 1
$ws.Range("AJ16").Value = 1
$ws.Range("AN16").Value = 0.06
$ws.Range("AO16").Value = 0.5
$ws.Range("AP16").Value = 0.5
$ws.Range("AQ16").Value = 0.04
$ws.Range("AR16").Value = 0.08
$ws.Range("AS16").Value = 0.04
$ws.Range("AT16").Value = 0.11
$ws.Range("AU16").Value = 0.04
$ws.Range("AV16").Value = 0.08
$ws.Range("AW16").Value = 0.7
$ws.Range("AX16").Value = 1
$ws.Range("AY16").Value = 15
$ws.Range("AZ16").Value = "closed"
$ws.Range("BA16").Value = "cd"
$ws.Range("BB16").Value = 0.02
$ws.Range("BC16").Value = "method1"
$ws.Range("BD16").Value = 5
$ws.Range("BE16").Value = 0.0822
$ws.Range("BF16").Value = 0.12
$ws.Range("BG16").Value = 0.075
$ws.Range("BH16").Value = 0.02
$ws.Range("BI16").Value = "constant"
$ws.Range("BJ16").Value = 0.02
$ws.Range("BK16").Value = 0.01
$ws.Range("BL16").Value = "AL_pct"
$ws.Range("BM16").Value = 0.75
$ws.Range("BN16").Value = $true
$ws.Range("BO16").Value = $true

# --- New row 17: DC_pure (name only) ---
$ws.Range("A17").Value = "DC_pure"

# --- Extend dropdown (list) data validation to cover the new row 16 ---
# H column: preSet / ALpct / MApct
$ws.Range("H16").Validation.Add(3, 1, 3, "preSet, ALpct,MApct")

# F column: list sourced from the Policies sheet
$ws.Range("F16").Validation.Add(3, 1, 3, "=Policies!`$B`$3:`$B`$8")

# --- Reflect the new selection / view position used while editing ---
$ws.Activate()
$ws.Range("C20").Select()
